$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the progress for the alinea in row 24: now at 10% and assigned to João.
$ws.Range("D24").Value2 = 10
$ws.Range("E24").Value2 = "João"

# Move the active selection from D27 to D25 as in the saved view state.
$ws.Range("D25").Select() | Out-Null

Write-Host "Applied edits to Progress worksheet"
